$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row-level updates derived from the published diff.
# Each entry maps row number -> column letter -> new text value.
$rowUpdates = @(
    @{ Row=2; D='43.111.75'; E='  -4.01%  ' },
    @{ Row=3; D='2.221.12'; E='  -5.36%  ' },
    @{ Row=4; E='  +0.05%  ' },
    @{ Row=5; D='317.25'; E='  -3.76%  ' },
    @{ Row=6; D='98.12'; E='  -7.25%  ' },
    @{ Row=7; D='0.579'; E='  -8.00%  ' },
    @{ Row=8; E='  +0.02%  ' },
    @{ Row=9; D='0.561'; E='  -7.79%  ' },
    @{ Row=10; D='36.87'; E='  -8.62%  ' },
    @{ Row=11; D='54.05'; E='  -3.13%  ' },
    @{ Row=12; E='  -9.52%  ' },
    @{ Row=13; D='7.65'; E='  -8.38%  ' },
    @{ Row=14; E='  -2.06%  ' },
    @{ Row=15; D='2.560.59'; E='  -5.29%  ' },
    @{ Row=16; D='0.859'; E='  -11.00%  ' },
    @{ Row=17; D='14.29'; E='  -6.17%  ' },
    @{ Row=18; D='2.227.20'; E='  -6.17%  ' },
    @{ Row=19; D='43.016.39'; E='  -4.18%  ' },
    @{ Row=20; D='13.68'; E='  -10.54%  ' },
    @{ Row=21; E='  -9.42%  ' },
    @{ Row=22; D='0.0₃0959'; E='  -9.12%  ' },
    @{ Row=23; E='  -11.12%  ' },
    @{ Row=24; D='65.02'; E='  -10.50%  ' },
    @{ Row=25; D='235.57'; E='  -8.40%  ' },
    @{ Row=26; D='2.18'; E='  -3.68%  ' },
    @{ Row=27; E='  +0.10%  ' },
    @{ Row=28; D='4.04'; E='  +1.54%  ' },
    @{ Row=29; D='10.00'; E='  -10.86%  ' },
    @{ Row=30; E='  -3.18%  ' },
    @{ Row=31; B='Filecoin'; C='https://coinranking.com/coin/ymQub4fuB+filecoin-fil'; D='6.41'; E='  -12.73%  ' },
    @{ Row=32; B='InjectiveProtocol'; C='https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'; D='36.64'; E='  -0.13%  ' },
    @{ Row=33; D='20.17'; E='  -8.00%  ' },
    @{ Row=34; D='0.0860'; E='  -9.55%  ' },
    @{ Row=35; D='157.21' },
    @{ Row=36; D='3.31'; E='  +2.53%  ' },
    @{ Row=37; D='2.67'; E='  -4.06%  ' },
    @{ Row=38; D='0.120'; E='  -8.12%  ' },
    @{ Row=39; D='1.84'; E='  -3.80%  ' },
    @{ Row=40; E='  -6.44%  ' },
    @{ Row=41; E='  -10.04%  ' },
    @{ Row=42; D='3.66'; E='  -6.64%  ' },
    @{ Row=43; D='0.0317'; E='  -9.20%  ' },
    @{ Row=44; D='13.92'; E='  +9.02%  ' },
    @{ Row=45; E='  +0.04%  ' },
    @{ Row=46; D='1.746.45'; E='  -6.96%  ' },
    @{ Row=47; B='FraxShare'; C='https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'; D='9.09'; E='  -2.19%  ' },
    @{ Row=48; B='Algorand'; C='https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'; D='0.201'; E='  -10.76%  ' },
    @{ Row=49; D='83.19'; E='  -12.53%  ' },
    @{ Row=50; D='5.25'; E='  -12.72%  ' },
    @{ Row=51; D='73.44'; E='  -12.38%  ' }
)

foreach ($u in $rowUpdates) {
    $r = $u.Row

    if ($u.ContainsKey("B")) {
        $ws.Range("B" + $r).Value = $u.B
    }
    if ($u.ContainsKey("C")) {
        $ws.Range("C" + $r).Value = $u.C
    }
    if ($u.ContainsKey("D")) {
        # Force text so numeric-looking prices (e.g. "317.25", "10.00")
        # stay stored as strings instead of being coerced to numbers.
        $dCell = $ws.Range("D" + $r)
        $dCell.NumberFormat = "@"
        $dCell.Value = $u.D
        $dCell.Style = "Normal"
    }
    if ($u.ContainsKey("E")) {
        $ws.Range("E" + $r).Value = $u.E
    }
}

Write-Output ("Updated " + $rowUpdates.Count + " rows")
